# Append two new "historico" rows (43 and 44) to the Historico sheet,
# matching the newsbot's latest scraped entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43 --------------------------------------------------------------
$ws.Range("A43").Value = "05/01/2026 06:03:05"
$ws.Range("B43").Value = "05/01 06:00"
$ws.Range("C43").Value = "Folha de S.Paulo - Poder - Principal"
$ws.Range("D43").Value = "Governo Lula reclama da Argentina em reunião sobre direitos humanos do Mercosul"
$ws.Range("E43").Value = "https://redir.folha.com.br/redir/online/poder/rss091/*https://www1.folha.uol.com.br/colunas/painel/2026/01/governo-lula-reclama-da-argentina-em-reuniao-sobre-direitos-humanos-do-mercosul.shtml"
$ws.Range("F43").Value = "lula"
$ws.Range("G43").Value = "os do Mercosul, realizada em outubro, gerou queixas do Ministério dos Direitos Humanos de Lula (PT).`n&lt;a href=`"https://redir.folha.com.br/redir/online/poder/rss091/*https://www1.folha.u"

# --- Row 44 ----------------------------------------------------------------
$ws.Range("A44").Value = "05/01/2026 06:03:06"
$ws.Range("B44").Value = "05/01 05:37"
$ws.Range("C44").Value = "BBC Brasil"
$ws.Range("D44").Value = "EUA 'governarem' a Venezuela não tem respaldo legal, a era das tutelas terminou com a descolonização, diz professora de Oxford"
$ws.Range("E44").Value = "https://www.bbc.com/portuguese/articles/c74vxlkm13zo?at_medium=RSS&at_campaign=rss"

# F44/G44 exist in the source as present-but-empty cells (no score/snippet
# extracted for this item yet). Touching a formatting property materializes
# the cell in the sheet without assigning it a value.
$ws.Range("F44").Font.Bold = $false
$ws.Range("G44").Font.Bold = $false
